$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'327.47"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'5.80%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'39.82"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'7.50%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.706"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'11.36%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08065"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.46%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'4.568"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'3.78%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.680"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'4.47%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'1.949"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'4.66%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'1.00%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9443"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.04%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1277"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'12.13%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1987"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'5.89%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09249"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'4.76%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03456"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'4.73%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.09616"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.46%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001310"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-5.25%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006126"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.35%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.376"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.55%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.3498"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1.39%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'7.631"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'19.24%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1412"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'9.29%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.2511"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'5.74%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04401"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'1.32%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001254"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'4.27%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004317"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.08%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001194"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-14.87%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0004002"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'37.76%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02516"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'17.95%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05204"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'4.75%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007314"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-3.22%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1427"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'5.59%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.009084"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'7.26%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002196"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'5.90%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.01006"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'25.81%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006746"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'2.45%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.14%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.002881"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-12.66%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.001806"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'24.91%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002107"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.14%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002006"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.14%"
$ws.Range("E51").Style = "Normal"
